$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '51.237.98'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.20%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.769.26'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.12%  '

$ws.Range("E4").Value = '  -0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '353.25'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.76%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '107.54'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.97%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.548'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -2.23%  '

$ws.Range("E8").Value = '  +0.00%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.581'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.14%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.41'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.61%  '

$ws.Range("E11").Value = '  +3.33%  '

$ws.Range("B12").Value = 'Dogecoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0831'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.16%  '

$ws.Range("B13").Value = 'Chainlink'
$ws.Range("C13").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '19.94'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +3.23%  '

$ws.Range("E14").Value = '  -0.94%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.205.70'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.28%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.764.38'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.88%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.923'
$ws.Range("D17").Style = "Normal"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '51.193.02'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.12%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.62'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.20%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.09'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.90%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.07'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.47%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.0₃0958'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.48%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '69.53'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.29%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '264.89'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -3.33%  '

$ws.Range("E25").Value = '  -0.55%  '

$ws.Range("E26").Value = '  +0.02%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '25.92'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.07%  '

$ws.Range("E28").Value = '  +13.28%  '

$ws.Range("B29").Value = 'Cosmos'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '10.14'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.34%  '

$ws.Range("B30").Value = 'Toncoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.24'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.04%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '35.54'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +5.84%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '51.76'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.78%  '

$ws.Range("E33").Value = '  +6.71%  '

$ws.Range("B34").Value = 'RenderToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.53'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +4.58%  '

$ws.Range("B35").Value = 'VeChain'
$ws.Range("C35").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0442'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.84%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0825'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.90%  '

$ws.Range("E37").Value = '  -0.13%  '

$ws.Range("E38").Value = '  +0.67%  '

$ws.Range("E39").Value = '  -2.02%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.96'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.52%  '

$ws.Range("E41").Value = '  -0.36%  '

$ws.Range("E42").Value = '  -1.03%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '120.92'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.33%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '21.96'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.62%  '

$ws.Range("E45").Value = '  -2.20%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.094.54'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.20%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.22'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.38%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.28'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.49%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.903'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.65%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '5.40'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -4.87%  '

$ws.Range("E51").Value = '  +7.37%  '
